$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date in column C for all existing data
# rows (2-116) from 45205 (2023-10-06) to 45206 (2023-10-07).
for ($r = 2; $r -le 116; $r++) {
    $ws.Cells.Item($r, 3).Value = 45206
}

# Row 116 gains an explicit row height (15pt, custom height).
$ws.Rows.Item(116).RowHeight = 15

# Add a new data row (117) for case "A 48309-2023".
# Copy formatting from row 116 first (split A:E / G:R so column F,
# which is never used in data rows, stays empty), then set the values.
$ws.Range("A116:E116").Copy()
$ws.Range("A117").PasteSpecial(-4122)
$ws.Range("G116:R116").Copy()
$ws.Range("G117").PasteSpecial(-4122)

$ws.Range("A117").Value = "A 48309-2023"
$ws.Range("B117").Value = 45205
$ws.Range("C117").Value = 45206
$ws.Range("D117").Value = "SKÅNE LÄN"
$ws.Range("E117").Value = "ÄNGELHOLM"
$ws.Range("G117").Value = 2.8
$ws.Range("H117").Value = 0
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 0
$ws.Range("N117").Value = 0
$ws.Range("O117").Value = 0
$ws.Range("P117").Value = 0
$ws.Range("Q117").Value = 0

$excel.CutCopyMode = 0
